# Weekly update: insert a new price record as row 123 in the "Camote"
# (Vega Modelo de Temuco) sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 123; Excel shifts rows 123:210 down
# to 124:211 and copies their formatting (incl. the date style on column D).
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(123, 1).Value  = 10
$ws.Cells.Item(123, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(123, 3).Value  = "La Araucanía"
$ws.Cells.Item(123, 4).Value  = 45086
$ws.Cells.Item(123, 5).Value  = 9
$ws.Cells.Item(123, 6).Value  = 100114002
$ws.Cells.Item(123, 7).Value  = "Camote"
$ws.Cells.Item(123, 8).Value  = "Sin especificar"
$ws.Cells.Item(123, 9).Value  = "Primera"
$ws.Cells.Item(123, 10).Value = 40
$ws.Cells.Item(123, 11).Value = 26000
$ws.Cells.Item(123, 12).Value = 26000
$ws.Cells.Item(123, 13).Value = 26000
$ws.Cells.Item(123, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(123, 15).Value = "Perú"
$ws.Cells.Item(123, 16).Value = 1444
$ws.Cells.Item(123, 17).Value = 18
$ws.Cells.Item(123, 18).Value = "Hortaliza"
